$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price/volume figures from the latest coinranking.com scrape.
# Values are stored as literal text (e.g. "282.49" / "1.54%"), so force
# text format on each touched cell before writing to avoid Excel
# reinterpreting the numeric-looking strings as numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "282.49"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.54%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.14%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.017"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.03%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06505"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.24%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.216"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.89%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.385"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "15.80%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9177"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.55%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1536"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.39%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06364"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "23.65%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07567"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.78%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02860"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.68%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08981"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.10%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001594"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.68%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006344"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.52%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006188"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.45%"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.92%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.359"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.63%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.242"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.16%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.07%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1326"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.13%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.987"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.10%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1543"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.61%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04442"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.87%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.91%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004439"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "14.40%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "1.71%"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.61%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04124"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.07%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006686"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.83%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1228"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.60%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002191"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "14.68%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01154"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.42%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005652"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.29%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.963"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "16.59%"
